$d = $word.ActiveDocument

# --- 1. Trim the "and it seems fine..." paragraph down to its first
#        sentence, merging the two runs into one. ---
$endash = [char]8211
$rsquo  = [char]8217
$hellip = [char]8230
$old = "and it seems fine. There is one thing left to do " + $endash + `
    " the maths library in the monitor ROM doesn" + $rsquo + `
    "t have < = > comparisons (in VTL-2 > is actually >= " + $hellip + $hellip + `
    ") so there" + $rsquo + "s a separate routine to do this which is, as yet, " + `
    "not implemented, or tested. Then I can start to think about the right expression. "
$new = "and it seems fine. "
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# --- 2/3. Move the "_GoBack" bookmark from the trailing empty paragraph
#          to the very start of the "Incidentally the :" paragraph.
#          Adding a bookmark with a name that already exists elsewhere
#          relocates it (bookmark names are unique), which removes it
#          from the final paragraph automatically. ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.StartsWith("Incidentally the")) {
        $target = $para
        break
    }
}
$r = $target.Range
$r.Collapse(1)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
